# Insert a new weekly price record into the daily logic sub-sheet.
# The new observation slots in right above the existing row 160, pushing
# every following row down by one (old row 227 becomes new row 228).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(160).Insert()

$ws.Range("A160").Value = 4
$ws.Range("B160").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C160").Value = "Los Lagos"
$ws.Range("D160").Value = 44609
$ws.Range("E160").Value = 10
$ws.Range("F160").Value = 100112043
$ws.Range("G160").Value = "Pepino ensalada"
$ws.Range("H160").Value = "Sin especificar"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 200
$ws.Range("K160").Value = 17000
$ws.Range("L160").Value = 17000
$ws.Range("M160").Value = 17000
$ws.Range("N160").Value = "$/caja 60 unidades"
$ws.Range("O160").Value = "Región de Arica y Parinacota"
$ws.Range("P160").Value = 283
$ws.Range("Q160").Value = 60
$ws.Range("R160").Value = "Hortaliza"
